# Update "want to go" counts (column F) for several rows across sheets.
# Values below mirror the refreshed scrape output committed to gh-pages.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3281
$ws1.Range("F4").Value = 2016
$ws1.Range("F6").Value = 119
$ws1.Range("F7").Value = 3101
$ws1.Range("F8").Value = 621
$ws1.Range("F12").Value = 156
$ws1.Range("F15").Value = 10244
$ws1.Range("F17").Value = 240
$ws1.Range("F18").Value = 16
$ws1.Range("F20").Value = 8132
$ws1.Range("F21").Value = 12746
$ws1.Range("F24").Value = 34
$ws1.Range("F27").Value = 603
$ws1.Range("F28").Value = 70
$ws1.Range("F30").Value = 2834
$ws1.Range("F33").Value = 8028
$ws1.Range("F34").Value = 1652
$ws1.Range("F37").Value = 89
$ws1.Range("F39").Value = 1488
$ws1.Range("F41").Value = 387
$ws1.Range("F42").Value = 86
$ws1.Range("F43").Value = 646

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 26

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3281
$ws4.Range("F6").Value = 2016
$ws4.Range("F9").Value = 26
$ws4.Range("F10").Value = 3101
$ws4.Range("F12").Value = 621
$ws4.Range("F15").Value = 156
$ws4.Range("F18").Value = 10244
$ws4.Range("F19").Value = 240
$ws4.Range("F20").Value = 16
$ws4.Range("F22").Value = 8132
$ws4.Range("F23").Value = 12746
$ws4.Range("F25").Value = 34
$ws4.Range("F28").Value = 603
$ws4.Range("F30").Value = 70
$ws4.Range("F31").Value = 2834
$ws4.Range("F36").Value = 8028
$ws4.Range("F39").Value = 89
$ws4.Range("F47").Value = 646
